# Add new data rows (18-22) to the "All channels" worksheet, mirroring the
# five extra subject-bins appended to the bin3_mean_amplitude_P3 export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used by every new row (A .. S == 19 channels).
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")

$newRows = @{
    18 = @(-0.69, 0.25, -0.4, -0.23, 0.96, 1.47, 0.45, -0.5, 1.19, -2.11, 0.05, -0.04, 0.43, -0.12, 0.08, 0, 0.01, -4.97, -0.32)
    19 = @(2.75, 2.08, -0.19, -0.02, -0.83, -0.55, 1.7, -0.63, -0.97, -1.35, -0.13, 0.31, 0.04, -0.25, -1.21, 0.12, 0.69, 1.01, 0.55)
    20 = @(-2.44, -2.2, -1, -1.73, -3.64, -1.15, 0.59, -1.38, -2, -1.17, 0.12, -2.12, -1.22, -0.9, -0.87, -1.9, -0.24, 0.77, 1.15)
    21 = @(2.5, -5.76, -1.39, 1.26, 2.38, -8.67, 0.68, 1.65, 0.48, 1.77, 1.44, 1.41, 2.41, 3.39, 1.5, 0.31, -1.92, -1.64, 0.79)
    22 = @(-1.36, -4.78, 1.09, 0.08, 4.97, 2.18, 1.07, 1.19, 0.65, 0.95, -0.89, -0.3, 2.06, 0.02, 2.57, 2.45, 0.3, 0.94, 1.05)
}

foreach ($r in 18, 19, 20, 21, 22) {
    $vals = $newRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
    # Newly-created cells otherwise pick up an explicit "applyFill" style
    # stamp; reset to Normal so they stay styleless like every other cell
    # on this sheet.
    $ws.Range("A" + $r + ":S" + $r).Style = "Normal"
}

# Move/collapse the selection to J19, matching the saved view state.
$ws.Range("J19").Select() | Out-Null
